{"js": "// ParserGrammar.docx edit \u2014 \"Declaration statments 80%, better error logging\"\n//\n// Semantic changes (paragraph text content; run/proofErr splitting left to\n// the editor and is not meaningful to reproduce by hand):\n//   1. <declaration statement> paragraph gains the \"foo; | foo := expr;\"\n//      alternation (declaring with and without an initializer).\n//   2. A previously-empty paragraph (between <float expression> and\n//      <relop>) now holds the new <Type_Name> production.\n//   3. <arithmetic op> gains a \"| %\" (modulo) alternative.\n//   4. The <bool op> line drops its dangling trailing \" | \".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// 1) <declaration statement> -> ... ; | <Var_Name> : <Type_Name> := <expression> ;\nitems[4].insertText(\n  \"<declaration statement> -> <Var_Name> : <Type_Name> ; | <Var_Name> : <Type_Name> := <expression> ;\",\n  Word.InsertLocation.replace\n);\n\n// 2) formerly-empty paragraph -> new <Type_Name> production\nitems[12].insertText(\n  \"<Type_Name> ->  Float | Integer | String | Character | Boolean\",\n  Word.InsertLocation.replace\n);\n\n// 3) <arithmetic op> -> ... | %\nitems[14].insertText(\" | %\", Word.InsertLocation.end);\n\n// 4) <bool op> -> drop trailing \" | \"\nitems[15].insertText(\n  \"!(ONLY MAYBE DOING THIS) <bool op> -> && | \\u201c||\\u201d\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# ParserGrammar.docx edit - \"Declaration statments 80%, better error logging\"\n#\n# Semantic changes (paragraph text content; run/proofErr splitting is\n# cosmetic markup the editor regenerates and is not reproduced by hand):\n#   1. <declaration statement> paragraph gains the \"foo; | foo := expr;\"\n#      alternation (declaring with and without an initializer).\n#   2. A previously-empty paragraph (between <float expression> and\n#      <relop>) now holds the new <Type_Name> production.\n#   3. <arithmetic op> gains a \"| %\" (modulo) alternative.\n#   4. The <bool op> line drops its dangling trailing \" | \".\n\n$d = $word.ActiveDocument\n\n# 1) <declaration statement> -> ... ; | <Var_Name> : <Type_Name> := <expression> ;\n$rDecl = $d.Paragraphs.Item(5).Range\n$rDecl.End = $rDecl.End - 1\n$rDecl.Text = \"<declaration statement> -> <Var_Name> : <Type_Name> ; | <Var_Name> : <Type_Name> := <expression> ;\"\n\n# 2) formerly-empty paragraph -> new <Type_Name> production\n$rType = $d.Paragraphs.Item(13).Range\n$rType.End = $rType.End - 1\n$rType.Text = \"<Type_Name> ->  Float | Integer | String | Character | Boolean\"\n\n# 3) <arithmetic op> -> ... | %\n$rArith = $d.Paragraphs.Item(15).Range\n$rArith.End = $rArith.End - 1\n$rArith.InsertAfter(\" | %\")\n\n# 4) <bool op> -> drop trailing \" | \"\n$leftCurly = [char]8220\n$rightCurly = [char]8221\n$newBoolOp = \"!(ONLY MAYBE DOING THIS) <bool op> -> && | \" + $leftCurly + \"||\" + $rightCurly\n$rBoolOp = $d.Paragraphs.Item(16).Range\n$rBoolOp.End = $rBoolOp.End - 1\n$rBoolOp.Text = $newBoolOp\n"}
